# ChipStatus.xlsx - Added SCP VIs and transfer queue (to be tested).
# Applies the worksheet edits to Feuil1 (sheet1): new highlight colors in
# column A, some text/content corrections in the "rejected by Promex" rows,
# a couple of new cells, and three new trailer rows (SLAC for board / SLAC
# for CSAM / CENBG).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlPasteFormats
$xlPasteFormats = -4122

# Excel VBA-style BGR color values for the two new highlight fills used in
# column A: orange (RGB 255,192,0) and blue (RGB 0,176,240).
$orange = 49407
$blue   = 15773696

# ---------------------------------------------------------------------
# 1) Column A highlight colors
# ---------------------------------------------------------------------
$orangeCells = @("A2","A63","A87","A88","A89","A90","A92","A93","A100","A105")
foreach ($addr in $orangeCells) {
    $ws.Range($addr).Interior.Color = $orange
}

$blueCells = @("A12","A13","A19","A23","A26","A46","A52","A55","A59","A65","A106")
foreach ($addr in $blueCells) {
    $ws.Range($addr).Interior.Color = $blue
}

# ---------------------------------------------------------------------
# 2) Column F "Rejected" formatting (re-use the existing red-on-green
#    style already used by F12/F13/F58/F100) for the newly / re-flagged
#    rejected rows.
# ---------------------------------------------------------------------
$ws.Range("F12").Copy()
$fFormatCells = @("F19","F23","F26","F46","F52","F55","F59")
foreach ($addr in $fFormatCells) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = 0

# F46, F55 and F59 are brand-new cells in column F -- give them the
# "Rejected" text (shared with D46/D55/D59).
$ws.Range("F46").Value = "Rejected"
$ws.Range("F55").Value = "Rejected"
$ws.Range("F59").Value = "Rejected"

# ---------------------------------------------------------------------
# 3) Column D "Rejected" formatting for the new D46 cell (re-use the
#    existing style already on D12/D13/D19/D50).
# ---------------------------------------------------------------------
$ws.Range("D12").Copy()
$ws.Range("D46").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("D46").Value = "Rejected"

# ---------------------------------------------------------------------
# 4) Text corrections / additions
# ---------------------------------------------------------------------
# Row 46 - new incident text + reason
$ws.Range("C46").Value = "rejected by Promex, bad channel (3 tests)"
$ws.Range("E46").Value = "same bad channel"

# Row 52 - reason updated from 3 tests to 4 tests
$ws.Range("E52").Value = "missing channel (4 tests)"

# Row 55 - new incident + reason
$ws.Range("C55").Value = "rejected by Promex, bad channel"
$ws.Range("E55").Value = "same bad channel"

# Row 59 - new incident + reason
$ws.Range("C59").Value = "rejected by Promex, no signal"
$ws.Range("E59").Value = "same no signal"

# F64 gets a single-space placeholder value (keeps its existing style).
$ws.Range("F64").Value = " "

# G88 / G90 - "OK 1,2,3" note (re-uses the existing shared text already in G2).
$ws.Range("G88").Value = "OK 1,2,3"
$ws.Range("G90").Value = "OK 1,2,3"

# ---------------------------------------------------------------------
# 5) New trailer rows: SLAC for board / SLAC for CSAM / CENBG
# ---------------------------------------------------------------------
# Row 104 re-uses the existing yellow "Pending" style already used on A3/A4.
$ws.Range("A3").Copy()
$ws.Range("A104").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("A104").Value = "SLAC for board"

$ws.Range("A105").Value = "SLAC for CSAM"
$ws.Range("A105").Interior.Color = $orange

$ws.Range("A106").Value = "CENBG"
$ws.Range("A106").Interior.Color = $blue

# ---------------------------------------------------------------------
# 6) View bookkeeping: selection / top-left cell reset to match the
#    saved workbook state.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A58:XFD58").Select()

Write-Output "edit complete"
